$d = $word.ActiveDocument

$replacements = @(
    @("74×74=", "93×48="),
    @("80×67=", "21×26="),
    @("83×18=", "93×24="),
    @("22×54=", "12×88="),
    @("15×35=", "72×24="),
    @("13×66=", "17×81="),
    @("22×74=", "29×99="),
    @("96×40=", "68×97="),
    @("44×62=", "28×39="),
    @("81×35=", "33×34="),
    @("80×93=", "77×96="),
    @("67×22=", "67×96="),
    @("81×62=", "19×54="),
    @("60×97=", "80×92="),
    @("50×52=", "60×37="),
    @("37×94=", "97×70="),
    @("69×53=", "43×95="),
    @("77×39=", "59×39="),
    @("37×36=", "83×13="),
    @("43×46=", "93×29="),
    @("39×65=", "84×29="),
    @("83×50=", "52×44="),
    @("89×41=", "80×89="),
    @("39×57=", "85×66="),
    @("78×54=", "59×51=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
